$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "unify the conception of DataNode, DataTable, Entity" -- the sheet that
# used to describe a single "Property1" table is renamed to "DataNode".
$ws.Name = "DataNode"

# Leave the cursor where the author was last working (D40) instead of the
# old A9 selection.
[void]$ws.Range("D40").Select()
